$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "K"  = -1.59
    "U"  = 1.07
    "V"  = 0.1903914590747331
    "W"  = -0.0001355533389600757
    "X"  = 2.400763686053281
    "Y"  = -2.400899239392241
    "AA" = -0.00004917896140701842
    "AB" = 0.3286911707262999
    "AC" = -0.3287403496877069
    "AD" = 82.8
    "AF" = 82.8
    "AG" = 81.73
    "AH" = 0.9364397195204704
    "AI" = 0.08113669769720724
    "AJ" = 0.9356611333714939
    "AK" = 0.08017225312184259
    "AL" = 6.94
    "AM" = 0.9100000000000001
    "AN" = -130.3937007874016
    "AO" = -0.09351585014409222
    "AP" = -128.7086614173228
    "AQ" = -0.7131868131868131
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
